$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff duplicates rows 11:13 (the James/Dana/Kevin block whose D column
# links to Jane Smith / John Doe / Samantha Black's photos) twice, appended
# as new rows 14:16 and then again as rows 17:19 -- i.e. the user selected
# rows 11:13, copied them, and pasted the block twice at the bottom.

$urlJaneSmith      = "https://raw.githubusercontent.com/mostafaalyCS/sheet/refs/heads/main/Jane%20Smith.jpg"
$urlJohnDoe        = "https://raw.githubusercontent.com/mostafaalyCS/sheet/refs/heads/main/John%20Doe.jpg"
$urlSamanthaBlack  = "https://raw.githubusercontent.com/mostafaalyCS/sheet/refs/heads/main/Samantha%20Black.jpg"

# source row -> destination row, row height, and the hyperlink target for
# that row's D cell. Two copies of the 11:13 block, back to back.
$pasteOrder = @(
    @{ Src = 11; Dst = 14; Height = 53; Url = $urlJaneSmith },
    @{ Src = 12; Dst = 15; Height = 53; Url = $urlJohnDoe },
    @{ Src = 13; Dst = 16; Height = 51; Url = $urlSamanthaBlack },
    @{ Src = 11; Dst = 17; Height = 53; Url = $urlJaneSmith },
    @{ Src = 12; Dst = 18; Height = 53; Url = $urlJohnDoe },
    @{ Src = 13; Dst = 19; Height = 51; Url = $urlSamanthaBlack }
)

foreach ($p in $pasteOrder) {
    [void]$ws.Rows.Item($p.Src).Copy()
    [void]$ws.Rows.Item($p.Dst).Insert(-4121)
    $ws.Rows.Item($p.Dst).RowHeight = $p.Height
}

# Hyperlinks.Add() resets the target cell's font to Excel's generic
# Hyperlink style, so re-paste just the source D cell's formatting
# afterwards to restore the original (themed, size-20) look -- matching
# how the rest of the D-column "image link" cells are styled.
# Order matches the rId order Excel itself assigned on paste: D16, D14,
# D15, D19, D17, D18.
$hyperlinkOrder = @(
    @{ Dst = "D16"; Src = "D13"; Url = $urlSamanthaBlack },
    @{ Dst = "D14"; Src = "D11"; Url = $urlJaneSmith },
    @{ Dst = "D15"; Src = "D12"; Url = $urlJohnDoe },
    @{ Dst = "D19"; Src = "D13"; Url = $urlSamanthaBlack },
    @{ Dst = "D17"; Src = "D11"; Url = $urlJaneSmith },
    @{ Dst = "D18"; Src = "D12"; Url = $urlJohnDoe }
)

foreach ($h in $hyperlinkOrder) {
    [void]$ws.Hyperlinks.Add($ws.Range($h.Dst), $h.Url)
    [void]$ws.Range($h.Src).Copy()
    [void]$ws.Range($h.Dst).PasteSpecial(-4122)
}

$ws.CutCopyMode = 0

# Mirror the author's final on-screen state: whole rows 17:19 (the
# just-pasted block) selected, with the window scrolled so row 6 is the
# top visible row.
[void]$ws.Range("A17:XFD19").Select()
$excel.ActiveWindow.ScrollRow = 6
